$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.635.24'
$ws.Range("E2").Value = '  +2.59%  '
$ws.Range("D3").Value = '2.700.60'
$ws.Range("E3").Value = '  +2.37%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'526.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").Value = "'145.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.65%  '
$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.81%  '
$ws.Range("D9").Value = '2.721.98'
$ws.Range("E9").Value = '  +2.52%  '
$ws.Range("D10").Value = "'6.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.65%  '
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("E13").Value = '  +3.09%  '
$ws.Range("D14").Value = '3.175.82'
$ws.Range("E14").Value = '  +2.30%  '
$ws.Range("D15").Value = '60.610.05'
$ws.Range("E15").Value = '  +2.57%  '
$ws.Range("D16").Value = "'21.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.16%  '
$ws.Range("D17").Value = '2.717.97'
$ws.Range("E17").Value = '  +2.73%  '
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("D19").Value = "'346.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.37%  '
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").Value = "'10.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.02%  '
$ws.Range("D22").Value = "'6.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.93%  '
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = "'63.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.40%  '
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("D28").Value = '0.0₃0821'
$ws.Range("E28").Value = '  +1.79%  '
$ws.Range("E29").Value = '  +2.08%  '
$ws.Range("D30").Value = "'6.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.99%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("E32").Value = '  +0.70%  '
$ws.Range("D33").Value = "'19.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").Value = "'150.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.26%  '
$ws.Range("E35").Value = '  +6.20%  '
$ws.Range("E36").Value = '  +8.46%  '
$ws.Range("D37").Value = "'0.938"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.31%  '
$ws.Range("D38").Value = "'0.875"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.12%  '
$ws.Range("E39").Value = '  +6.76%  '
$ws.Range("D40").Value = "'37.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.16%  '
$ws.Range("E41").Value = '  -0.58%  '
$ws.Range("D42").Value = "'283.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.39%  '
$ws.Range("D43").Value = "'20.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.50%  '
$ws.Range("D44").Value = '2.147.31'
$ws.Range("E44").Value = '  +7.77%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = "'0.611"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").Value = "'0.996"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").Value = "'0.0986"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("E48").Value = '  +2.31%  '
$ws.Range("D49").Value = "'4.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.02%  '
$ws.Range("D50").Value = "'10.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.67%  '
$ws.Range("E51").Value = '  +0.86%  '
